$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# Row 1 - headers
# Old layout: A Judges Names | B Appeleant Names | C Year | D Outcome Of Case
#             | E Type of Issue | F Gender of Appellant | G Gender of Judge
# New layout: ... | F Inital Comparison | G Duration of the Case
#             | H Type of Tax Payer | I Gender of Appellant | J Gender of Judge
#             | K URL
# ----------------------------------------------------------------------
$ws.Range("I1").Value = $ws.Range("F1").Value()
$ws.Range("J1").Value = $ws.Range("G1").Value()
$ws.Range("F1").Value = "Inital Comparison"
$ws.Range("G1").Value = "Duration of the Case"
$ws.Range("H1").Value = "Type of Tax Payer"
$ws.Range("K1").Value = "URL"

# ----------------------------------------------------------------------
# Row 2 - Lucie Lamarre / MANON RODIER (A-E unchanged)
# ----------------------------------------------------------------------
$ws.Range("I2").Value = $ws.Range("F2").Value()
$ws.Range("J2").Value = $ws.Range("G2").Value()
$ws.Range("F2").Value = 6
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "1"
$ws.Range("H2").Value = "Individual"
$ws.Range("K2").Formula = "=https://decision.tcc-cci.gc.ca/tcc-cci/decisions/en/item/26605/index.do"

# ----------------------------------------------------------------------
# Row 3 - previously " D.W. Beaubier" / CHARLES W. DOERING -> now
#          D.G.H. Bowman / CHARLES B. LOEWEN
# ----------------------------------------------------------------------
$ws.Range("A3").Value = "D.G.H. Bowman"
$ws.Range("B3").Value = "CHARLES B. LOEWEN,"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "2003"
$ws.Range("D3").Value = "Partially Winning"
$ws.Range("E3").Value = "Income Tax"
$ws.Range("F3").Value = 6
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "2"
$ws.Range("H3").Value = "Individual"
$ws.Range("I3").Value = "Male"
$ws.Range("J3").Value = "Male"

# ----------------------------------------------------------------------
# Row 4 - previously "C.H. McArthur" / MICHAEL LAURIE -> now
#          " D.W. Beaubier" / DOUGLAS DIXON,
# ----------------------------------------------------------------------
$ws.Range("A4").Value = " D.W. Beaubier"
$ws.Range("B4").Value = "DOUGLAS DIXON,"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "2003"
$ws.Range("D4").Value = "Partially Winning"
$ws.Range("E4").Value = "Income Tax"
$ws.Range("F4").Value = 2
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "Less than 1 year"
$ws.Range("H4").Value = "Partnership"
$ws.Range("I4").Value = "Male"
$ws.Range("J4").Value = "Male"

# ----------------------------------------------------------------------
# Row 5 - previously " D.W. Beaubier" / DOUGLAS DIXON -> now
#          " Diane Campbell" / S.K. MANAGEMENT INC.,&&&#CORPORATION
# ----------------------------------------------------------------------
$ws.Range("A5").Value = " Diane Campbell"
$ws.Range("B5").Value = "S.K. MANAGEMENT INC.,&&&#CORPORATION"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "2003"
$ws.Range("D5").Value = "Losing"
$ws.Range("E5").Value = "Excise tax"
$ws.Range("F5").Value = 6
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "Less than 1 year"
$ws.Range("H5").Value = "Corporation"
$ws.Range("I5").Value = "Male"
$ws.Range("J5").Value = "Female"

# ----------------------------------------------------------------------
# Row 6 - removed entirely (A.A. Sarchuk / HONORA ZAKRISON)
# ----------------------------------------------------------------------
$ws.Rows("6:6").Delete()
